$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the existing "_GoBack" bookmark. In the original document it
#    sits at the end of the "Re-open the Console and File pop-outs."
#    paragraph; the edit moves it onto the "Open the ... file" bullet
#    instead, so the old one must go before we re-add it elsewhere.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Find the bullet that currently reads "Open the NEWFILE file" and
#    rewrite its runs so the filename becomes "testFile" (flagged with
#    spellStart/spellEnd proofing marks, matching how the document
#    marks other camelCase tokens such as "usabilityTest"). Re-add the
#    "_GoBack" bookmark right after the new text, at the end of the
#    paragraph's run content.
# ------------------------------------------------------------------
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Open the NEWFILE file*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Open the NEWFILE file' paragraph"
}

$full = $target.Range
# Exclude the trailing paragraph mark so the paragraph properties
# (list style/numbering) are left untouched by the replacement.
$body = $d.Range($full.Start, $full.End - 1)

$xml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Open the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t>testFile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> file</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$body.InsertXML($xml)

Write-Output ("Updated paragraph text: " + $target.Range.Text)
Write-Output ("_GoBack exists: " + $d.Bookmarks.Exists("_GoBack"))
